# Documento Inicial Herramientas y Tecnologias
#
# Insert a new document "Herramientas y Tecnologias" into column B at row 6,
# pushing the existing entries (Plan de Riesgo, Plan de Calidad, Plan de Pruebas)
# down one row (rows 7, 8, 9 respectively).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing values in B6:B8 down to B7:B9 (bottom-up so we don't
# overwrite a value before it has been copied further down).
$ws.Range("B9").Value = $ws.Range("B8").Value2
$ws.Range("B8").Value = $ws.Range("B7").Value2
$ws.Range("B7").Value = $ws.Range("B6").Value2

# Write the new entry into the now-vacated row.
$ws.Range("B6").Value = "Herramientas y Tecnologias"

# Restore the active cell selection as recorded in the edited workbook.
$ws.Range("G10").Select()
